$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-11 with the new values (row 1 header and row 2 stay unchanged)
$ws.Range("A3").Value = "-"
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "العراق "

$ws.Range("A4").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "لبنان "

$ws.Range("A5").Value = "الاستخبارات والمراقبة والرصد (مترجمه)"
$ws.Range("B5").Value = "بلد "
$ws.Range("C5").Value = "نهاريا "

$ws.Range("A6").Value = "لبنان "
$ws.Range("B6").Value = "بلد "
$ws.Range("C6").Value = "صبرة (مترجمه)"

$ws.Range("A7").Value = "الانتداب البريطاني على فلسطين "
$ws.Range("B7").Value = "علم الوجود الجغرافي السياسي "
$ws.Range("C7").Value = "صلحا "

$ws.Range("A8").Value = "-"
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "سنغافورة "

$ws.Range("A9").Value = "-"
$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = "الاحتلال الإسرائيلي لجنوب لبنان "

$ws.Range("A10").Value = "-"
$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = "جنوب لبنان "

$ws.Range("A11").Value = "لبنان "
$ws.Range("B11").Value = "بلد "
$ws.Range("C11").Value = "بيروت "

# Remove rows 12-14, which are no longer present in the updated sheet
$ws.Range("A12:C14").EntireRow.Delete()
